$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 ("1,500" -> "1500") must stay a text value (not get auto-converted to a
# number by Excel's type inference). Temporarily force a text number format,
# assign the value, then restore the cell's original style so no visible
# formatting change is left behind.
$origStyleB2 = $ws.Range("B2").Style
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1500"
$ws.Range("B2").Style = $origStyleB2

$ws.Range("C2").Value = "No, Soleo Health does not encompass community sites. Soleo Health is a specialty infusion pharmacy that focuses on providing complex pharmaceutical care in the home or at an alternate site of care, rather than community-based services."

$ws.Range("D2").Value = "No, Soleo Health is not influential on state or local policy. The society's primary focus is on providing healthcare services, rather than lobbying or advocacy efforts that influence policy decisions."

$ws.Range("E2").Value = "No, Soleo Health does not provide engagement opportunity with leadership. The company's leadership may not have a direct engagement opportunity with the society's members."

$ws.Range("F2").Value = "No, Soleo Health does not provide support for clinical trial recruitment. Soleo Health is a specialty pharmacy that focuses on providing infusion therapy services rather than clinical trial recruitment."

$ws.Range("G2").Value = "No, Soleo Health does not provide engagement opportunities with payors. They primarily focus on providing specialized pharmacy services to patients."

$ws.Range("H2").Value = "No, justification: Not specified in public information."

$ws.Range("I2").Value = "No, Soleo Health is not involved in therapeutic research collaborations. Soleo Health primarily focuses on providing home and alternate-site infusion services to patients."

$ws.Range("J2").Value = "No, justification: Soleo Health does not publicly disclose information on the expertise of its board members."

$ws.Range("L2").Value = "2025-03-17 06:56:33"
